$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fruta / Hortaliza market data.
# Update existing rows (2-22) whose Fecha / price / origin data changed.

# Row 2
$ws.Range("D2").Value = 44188
$ws.Range("K2").Value = 38000
$ws.Range("L2").Value = 40000
$ws.Range("M2").Value = 39000
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 1560

# Row 3
$ws.Range("D3").Value = 44160
$ws.Range("K3").Value = 28000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = 29000
$ws.Range("N3").Value = '$/malla 25 kilos'
$ws.Range("O3").Value = 'Región de O''Higgins'
$ws.Range("P3").Value = 1160

# Row 4
$ws.Range("D4").Value = 44323
$ws.Range("H4").Value = 'Magnum'
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21000
$ws.Range("N4").Value = '$/malla 25 kilos'
$ws.Range("O4").Value = 'Perú'
$ws.Range("P4").Value = 840

# Row 5
$ws.Range("D5").Value = 44363
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 26000
$ws.Range("M5").Value = 25500
$ws.Range("N5").Value = '$/malla 25 kilos'
$ws.Range("O5").Value = 'Perú'
$ws.Range("P5").Value = 1020

# Row 6
$ws.Range("D6").Value = 44203
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 21000
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 840

# Row 7
$ws.Range("D7").Value = 44230
$ws.Range("H7").Value = 'Magnum'
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 23000
$ws.Range("N7").Value = '$/saco 25 kilos'
$ws.Range("O7").Value = 'Región del Maule'
$ws.Range("P7").Value = 920

# Row 8
$ws.Range("D8").Value = 44265
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = 'Región Metropolitana'

# Row 9
$ws.Range("D9").Value = 44384
$ws.Range("H9").Value = 'Sin especificar'

# Row 10
$ws.Range("D10").Value = 44272
$ws.Range("H10").Value = 'Magnum'
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 24000
$ws.Range("M10").Value = 23000
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("P10").Value = 920

# Row 11
$ws.Range("D11").Value = 44433
$ws.Range("K11").Value = 25000
$ws.Range("L11").Value = 26000
$ws.Range("M11").Value = 25500
$ws.Range("N11").Value = '$/malla 25 kilos'
$ws.Range("O11").Value = 'Perú'
$ws.Range("P11").Value = 1020

# Row 12
$ws.Range("D12").Value = 44321
$ws.Range("M12").Value = 24500
$ws.Range("N12").Value = '$/saco 25 kilos'
$ws.Range("O12").Value = 'Región del Maule'
$ws.Range("P12").Value = 980

# Row 13
$ws.Range("D13").Value = 44335
$ws.Range("K13").Value = 35000
$ws.Range("L13").Value = 36000
$ws.Range("M13").Value = 35500
$ws.Range("N13").Value = '$/saco 25 kilos'
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("P13").Value = 1420

# Row 14
$ws.Range("D14").Value = 44253
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 26000
$ws.Range("M14").Value = 25500
$ws.Range("P14").Value = 1020

# Row 15
$ws.Range("D15").Value = 44244
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17000
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 680

# Row 16
$ws.Range("D16").Value = 44441
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 29000
$ws.Range("M16").Value = 28500
$ws.Range("P16").Value = 1140

# Row 17
$ws.Range("D17").Value = 44237
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21000
$ws.Range("O17").Value = 'Región del Maule'
$ws.Range("P17").Value = 840

# Row 18
$ws.Range("D18").Value = 44342
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 28000
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = 29000
$ws.Range("N18").Value = '$/malla 25 kilos'
$ws.Range("O18").Value = 'Región Metropolitana'
$ws.Range("P18").Value = 1160

# Row 19
$ws.Range("D19").Value = 44294
$ws.Range("K19").Value = 24000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 24500
$ws.Range("N19").Value = '$/saco 25 kilos'
$ws.Range("O19").Value = 'Región del Maule'
$ws.Range("P19").Value = 980

# Row 20
$ws.Range("D20").Value = 44435
$ws.Range("K20").Value = 25000
$ws.Range("L20").Value = 26000
$ws.Range("M20").Value = 25500
$ws.Range("N20").Value = '$/malla 25 kilos'
$ws.Range("O20").Value = 'Perú'
$ws.Range("P20").Value = 1020

# Row 21
$ws.Range("D21").Value = 44279
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160

# Row 22
$ws.Range("D22").Value = 44167
$ws.Range("H22").Value = 'Sin especificar'
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 19000
$ws.Range("M22").Value = 18500
$ws.Range("P22").Value = 740

# Append new record (row 23) for the latest week.
$newRow = 23
$ws.Range("A$newRow").Value = 11
$ws.Range("B$newRow").Value = 'Vega Monumental Concepción'
$ws.Range("C$newRow").Value = 'Bíobío'
$ws.Range("D$newRow").Value = 44399
$ws.Range("D$newRow").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("E$newRow").Value = 8
$ws.Range("F$newRow").Value = 100112031
$ws.Range("G$newRow").Value = 'Poroto verde'
$ws.Range("H$newRow").Value = 'Magnum'
$ws.Range("I$newRow").Value = 'Primera'
$ws.Range("J$newRow").Value = 100
$ws.Range("K$newRow").Value = 20000
$ws.Range("L$newRow").Value = 22000
$ws.Range("M$newRow").Value = 21000
$ws.Range("N$newRow").Value = '$/malla 25 kilos'
$ws.Range("O$newRow").Value = 'Perú'
$ws.Range("P$newRow").Value = 840
$ws.Range("Q$newRow").Value = 25
$ws.Range("R$newRow").Value = 'Hortaliza'

Write-Host "Applied weekly Fruta/Hortaliza update"